# Spine RTS-GMLC output_parameters.xlsx update
# Renames the output-variable object names (obj_output / rel_report__output)
# to match the new Spine model's parameter names, and adds a new
# "node_state" output row. Also moves the active sheet from
# rel_report__output back to obj_output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("obj_output")
$ws3 = $wb.Worksheets.Item("rel_report__output")

# --- obj_output (sheet1) --------------------------------------------------
# Existing 7 rows get renamed in place; an 8th row (node_state) is new.
# B2 (was "flow") and the new B8 (node_state) keep a quote-prefixed style;
# B5 (units_started_up) keeps its quote-prefixed "d-mmm" number format.
# "node_state" is written first so new shared-string entries land in the
# same order the original author's Excel produced them.

$ws1.Cells.Item(8,1).Value = "output"
$ws1.Cells.Item(8,2).NumberFormat = "d-mmm"
$ws1.Cells.Item(8,2).Value = "'node_state"

$ws1.Cells.Item(2,1).Value = "output"
$ws1.Cells.Item(2,2).Value = "'unit_flow"

$ws1.Cells.Item(3,1).Value = "output"
$ws1.Cells.Item(3,2).ClearFormats()
$ws1.Cells.Item(3,2).Value = "connection_flow"

$ws1.Cells.Item(4,1).Value = "output"
$ws1.Cells.Item(4,2).Value = "units_on"

$ws1.Cells.Item(5,1).Value = "output"
$ws1.Cells.Item(5,2).NumberFormat = "d-mmm"
$ws1.Cells.Item(5,2).Value = "'units_started_up"

$ws1.Cells.Item(6,1).Value = "output"
$ws1.Cells.Item(6,2).Value = "units_shut_down"

$ws1.Cells.Item(7,1).Value = "output"
$ws1.Cells.Item(7,2).ClearFormats()
$ws1.Cells.Item(7,2).Value = "units_available"

# --- rel_report__output (sheet3) ------------------------------------------
# Insert a new row 8 (pushes the filler/formatting rows below down by one,
# preserving their style banding) then fill in the renamed C-column values
# the same way as sheet1, plus the new row's A/B/C values.

$ws3.Rows.Item(8).Insert()

$ws3.Cells.Item(2,3).Value = "'unit_flow"

$ws3.Cells.Item(3,3).ClearFormats()
$ws3.Cells.Item(3,3).Value = "connection_flow"

$ws3.Cells.Item(4,3).Value = "units_on"

$ws3.Cells.Item(5,3).NumberFormat = "d-mmm"
$ws3.Cells.Item(5,3).Value = "'units_started_up"

$ws3.Cells.Item(6,3).Value = "units_shut_down"

$ws3.Cells.Item(7,3).ClearFormats()
$ws3.Cells.Item(7,3).Value = "units_available"

$ws3.Cells.Item(8,1).Value = "report__output"
$ws3.Cells.Item(8,2).Value = "'report_1"
$ws3.Cells.Item(8,3).Value = "'node_state"

# --- view state -------------------------------------------------------
# Active tab moves from rel_report__output back to obj_output; selections
# move to the newly added rows/columns on each sheet.

$ws3.Range("C2:C8").Select()
$ws1.Activate()
$ws1.Range("A7:A8").Select()
